# @williamjreid Updated E5 from notes, added another results table
$wb = $excel.ActiveWorkbook

# --- Rename the first two sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Exercise 2"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Exercise 4"

# --- Populate "Exercise 4" (was empty Sheet2) with three small results tables ---

# Table 1
$ws2.Range("A1").Value = "a280_279"
$ws2.Range("A2").Value = 1667.1041398136299
$ws2.Range("A3").Value = 4127.2632729134702
$ws2.Range("A4").Value = 3869.2032729134698
$ws2.Range("A5").Value = 1235.46984738984
$ws2.Range("A6").Value = 231.27984738984699
$ws2.Range("A7").Value = 4688.2632729134702
$ws2.Range("A8").Value = 4205.8032729134702
$ws2.Range("A9").Value = 489.00413981362999
$ws2.Range("A10").Value = 388.35984738984899
$ws2.Range("A11").Value = 2462.83998930708
$ws2.Range("A13").Formula = "=AVERAGE(A2:A11)"
$ws2.Range("B13").Formula = "=MAX(A2:A11)"

# Table 2
$ws2.Range("A15").Value = "a280_1395"
$ws2.Range("A16").Value = -264555.3
$ws2.Range("A17").Value = -257430.7
$ws2.Range("A18").Value = -261865.4
$ws2.Range("A19").Value = -272915.8
$ws2.Range("A20").Value = -260774.9
$ws2.Range("A21").Value = -248197.8
$ws2.Range("A22").Value = -254304.6
$ws2.Range("A23").Value = -264046.40000000002
$ws2.Range("A24").Value = -263101.3
$ws2.Range("A25").Value = -262301.59999999998
$ws2.Range("A27").Formula = "=AVERAGE(A16:A25)"
$ws2.Range("B27").Formula = "=MAX(A16:A25)"

# Table 3
$ws2.Range("A29").Value = "a280_2790"
$ws2.Range("A30").Value = -731106.18
$ws2.Range("A31").Value = -747580.05
$ws2.Range("A32").Value = -787200.75
$ws2.Range("A33").Value = -737779.14
$ws2.Range("A34").Value = -736110.9
$ws2.Range("A35").Value = -722973.51
$ws2.Range("A36").Value = -741949.74
$ws2.Range("A37").Value = -712129.95
$ws2.Range("A38").Value = -727978.23
$ws2.Range("A39").Value = -714215.25
$ws2.Range("A41").Formula = "=AVERAGE(A30:A39)"
$ws2.Range("B41").Formula = "=MAX(A30:A39)"

# Column widths to match Exercise 2's layout
$ws2.Range("A1:B1").ColumnWidth = 13.8333333

# --- View / selection state ---
# Exercise 2: scroll window so row 4 is the top visible row, and select A23:B23
# (Exercise 2 is no longer the active tab once Exercise 4 is selected below)
$ws1.Activate()
$ws1.Range("A23:B23").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Exercise 4 becomes the active tab, with A41 selected
$ws2.Activate()
$ws2.Range("A41").Select()
